$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: date slash -> dash, plus D3 0->1, G3 0->1 ---
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# --- Row 4: date slash -> dash, plus D4 0->1, E4 0->1, H4 1->0 ---
# "01-08-2022" is ambiguous as a date (day<=12), so Excel would silently
# reinterpret it as a date serial; force the cell to Text first so the
# dash-formatted string is kept as a literal string, matching the source.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# --- Row 5: date slash -> dash, plus D5 0->1, E5 0->1, H5 1->0 ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# --- Rows 6-21: only the date format changes (slash -> dash); counts unchanged ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08-08-2022"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "11-08-2022"

$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "01-09-2022"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "05-09-2022"

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "08-09-2022"

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "12-09-2022"

$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"
